# Fix the invalid field names
#
# The "Sample Section" worksheet had a spurious, hidden, leading
# "header_info" column (A). This column is removed, which shifts every
# other column one position to the left (B->A, C->B, ... V->U).
#
# Column deletion in this COM-interop runtime correctly shifts cell
# *values* but does not relocate the cell *comments*, so the comments
# need to be captured and re-applied by hand at their new locations.
#
# In addition, the `.metadata` sheet's pav:createdOn timestamp is bumped
# to reflect the date this fix was made.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sample Section")

$oldCols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")
$newCols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U")

# 1. Capture the text of every comment in row 1 while columns are still in
#    their original positions, then delete the comments so nothing is left
#    dangling once the column shift happens.
$commentTexts = @{}
foreach ($col in $oldCols) {
    $cell = $ws.Range($col + "1")
    if ($cell.Comment -ne $null) {
        $commentTexts[$col] = $cell.Comment.Text()
        $cell.Comment.Delete()
    }
}

# 2. Delete column A outright (the hidden, spurious "header_info" column).
#    This shifts every column B..V left by one (B->A, C->B, ..., V->U),
#    including cell values, styles, widths and data validations.
$ws.Columns.Item(1).Delete()

# 3. Re-create the comments at their new (shifted) locations. The comment
#    that used to live on column A (blank placeholder) is intentionally
#    dropped, since that column itself no longer exists.
for ($i = 0; $i -lt $newCols.Count; $i++) {
    $oldCol = $oldCols[$i + 1]
    $newCol = $newCols[$i]
    if ($commentTexts.ContainsKey($oldCol)) {
        $ws.Range($newCol + "1").AddComment($commentTexts[$oldCol])
    }
}

# 4. Bump the pav:createdOn timestamp recorded on the `.metadata` sheet.
$meta = $wb.Worksheets.Item(".metadata")
$meta.Range("C2").Value = "2023-10-03T09:51:12-07:00"
